$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = -1168658269.12
$ws.Range("P2").Value = 3814468574.03
$ws.Range("Q2").Value = 3719126624
$ws.Range("R2").Value = -0.7394631714
$ws.Range("S2").Value = 2147242828.95
$ws.Range("T2").Value = 2147242828.95
$ws.Range("U2").Value = -0.2496171343
$ws.Range("V2").Value = 651267016.92
$ws.Range("W2").Value = 528325265.68
$ws.Range("X2").Value = 328884229.23
$ws.Range("Y2").Value = -935384416.33
$ws.Range("Z2").Value = -876133105.46
$ws.Range("AA2").Value = 108389636.21
$ws.Range("AG2").Value = 15450349.44
$ws.Range("AP2").Value = -3.3400478541
$ws.Range("AQ2").Value = -235.737917376023
$ws.Range("AR2").Value = -241.093925287188
$ws.Range("AS2").Value = -1181283529.76
$ws.Range("AT2").Value = -1001.39869006647
